$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (scenario 1, "Mean of absolute % of difference")
$ws.Range("C2").Value = 65
$ws.Range("D2").Value = 36
$ws.Range("F2").Value = 29
$ws.Range("G2").Value = 29
$ws.Range("H2").Value = 31
$ws.Range("I2").Value = 39

# Row 5 (scenario 2, "Mean of absolute % of difference")
$ws.Range("C5").Value = 34
$ws.Range("D5").Value = 37
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 38

# Row 6 (scenario 2, "min")
$ws.Range("C6").Value = -46
$ws.Range("D6").Value = -85
$ws.Range("E6").Value = -14
$ws.Range("F6").Value = -24
$ws.Range("G6").Value = -69
$ws.Range("H6").Value = -1
$ws.Range("I6").Value = -89

# Row 7 (scenario 2, "max")
$ws.Range("C7").Value = 79
$ws.Range("D7").Value = 51
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = 42
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 53

# Row 8 (scenario 3, "Mean of absolute % of difference")
$ws.Range("C8").Value = 23
$ws.Range("D8").Value = 30
$ws.Range("E8").Value = 18
$ws.Range("F8").Value = 23
$ws.Range("G8").Value = 25
$ws.Range("H8").Value = 50
$ws.Range("I8").Value = 30

# Row 9 (scenario 3, "min")
$ws.Range("C9").Value = -71
$ws.Range("D9").Value = -95
$ws.Range("E9").Value = -57
$ws.Range("F9").Value = -9
$ws.Range("G9").Value = -79
$ws.Range("H9").Value = -19
$ws.Range("I9").Value = -93

# Row 10 (scenario 3, "max")
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 71
$ws.Range("G10").Value = 11
$ws.Range("H10").Value = 157
$ws.Range("I10").Value = 12

# Remove rows 11-13 (scenario 4 block entirely removed)
$ws.Range("A11:I13").EntireRow.Delete()
